# Publish terminology IG 2.0.2 — update the ValueSet metadata table
# (sheet "Metadata") to match the new release: version, status, date,
# and clear the now-unused "Experimental" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (blank)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# (leading apostrophe forces literal text so the ISO-looking date string
# isn't auto-converted into a date serial number, matching the original
# plain-text cell content)
$ws.Range("B8").Value = "'2025-11-18"
